# Update loading_percent values for the 380 kV case (Case_4_213)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 10.71451415296707
$ws.Cells.Item(2, 3).Value = 9.502539268208498
$ws.Cells.Item(2, 4).Value = 5.99765845886862
$ws.Cells.Item(2, 5).Value = 12.82782147016723
$ws.Cells.Item(2, 7).Value = 33.85433203422809
$ws.Cells.Item(2, 8).Value = 15.66675799458318
$ws.Cells.Item(2, 9).Value = 24.99322015531614
$ws.Cells.Item(2, 11).Value = 8.117205471385432
$ws.Cells.Item(2, 12).Value = 10.14611150296172
$ws.Cells.Item(2, 13).Value = 13.40096278435721
$ws.Cells.Item(2, 15).Value = 24.50704399278917
$ws.Cells.Item(3, 2).Value = 10.45107117536305
$ws.Cells.Item(3, 3).Value = 9.47744563642966
$ws.Cells.Item(3, 4).Value = 5.88082069906799
$ws.Cells.Item(3, 5).Value = 12.86047377233772
$ws.Cells.Item(3, 7).Value = 33.95102876744952
$ws.Cells.Item(3, 8).Value = 15.71559616277484
$ws.Cells.Item(3, 9).Value = 25.0887774979617
$ws.Cells.Item(3, 11).Value = 7.913554400416784
$ws.Cells.Item(3, 12).Value = 10.15396045424003
$ws.Cells.Item(3, 13).Value = 13.36068383434054
$ws.Cells.Item(3, 15).Value = 24.58836856522416
$ws.Cells.Item(4, 2).Value = 10.28736439765946
$ws.Cells.Item(4, 3).Value = 9.462003878236246
$ws.Cells.Item(4, 4).Value = 5.80962127721272
$ws.Cells.Item(4, 5).Value = 12.88223737377831
$ws.Cells.Item(4, 7).Value = 34.01966845525153
$ws.Cells.Item(4, 8).Value = 15.74786791623343
$ws.Cells.Item(4, 9).Value = 25.15172439821565
$ws.Cells.Item(4, 11).Value = 7.786700132948081
$ws.Cells.Item(4, 12).Value = 10.16015981619323
$ws.Cells.Item(4, 13).Value = 13.33766270701587
$ws.Cells.Item(4, 15).Value = 24.64298669115745
$ws.Cells.Item(5, 2).Value = 10.22025858568139
$ws.Cells.Item(5, 3).Value = 9.455704253354872
$ws.Cells.Item(5, 4).Value = 5.780787517703414
$ws.Cells.Item(5, 5).Value = 12.89153778874193
$ws.Cells.Item(5, 7).Value = 34.04996256265508
$ws.Cells.Item(5, 8).Value = 15.76159370087677
$ws.Cells.Item(5, 9).Value = 25.17845050126894
$ws.Cells.Item(5, 11).Value = 7.734623473486151
$ws.Cells.Item(5, 12).Value = 10.16303370241941
$ws.Cells.Item(5, 13).Value = 13.32871789734759
$ws.Cells.Item(5, 15).Value = 24.66642064083018
$ws.Cells.Item(6, 2).Value = 10.20909483884275
$ws.Cells.Item(6, 3).Value = 9.454657830331282
$ws.Cells.Item(6, 4).Value = 5.776011921629903
$ws.Cells.Item(6, 5).Value = 12.89310819096544
$ws.Cells.Item(6, 7).Value = 34.05513294165553
$ws.Cells.Item(6, 8).Value = 15.76390757902925
$ws.Cells.Item(6, 9).Value = 25.18295325949989
$ws.Cells.Item(6, 11).Value = 7.72595528815391
$ws.Cells.Item(6, 12).Value = 10.16353191815374
$ws.Cells.Item(6, 13).Value = 13.32725916151075
$ws.Cells.Item(6, 15).Value = 24.67038285533577
$ws.Cells.Item(7, 2).Value = 10.28646084915843
$ws.Cells.Item(7, 3).Value = 9.461918945454183
$ws.Cells.Item(7, 4).Value = 5.809231624190124
$ws.Cells.Item(7, 5).Value = 12.88236105457409
$ws.Cells.Item(7, 7).Value = 34.02006761729174
$ws.Cells.Item(7, 8).Value = 15.74805069939785
$ws.Cells.Item(7, 9).Value = 25.15208048434005
$ws.Cells.Item(7, 11).Value = 7.785999258447772
$ws.Cells.Item(7, 12).Value = 10.16019716632636
$ws.Cells.Item(7, 13).Value = 13.33754029866327
$ws.Cells.Item(7, 15).Value = 24.64329796731559
$ws.Cells.Item(8, 2).Value = 10.62414073763918
$ws.Cells.Item(8, 3).Value = 9.493895041430257
$ws.Cells.Item(8, 4).Value = 5.95728764430841
$ws.Cells.Item(8, 5).Value = 12.83872424648675
$ws.Cells.Item(8, 7).Value = 33.88574524791883
$ws.Cells.Item(8, 8).Value = 15.68312329369101
$ws.Cells.Item(8, 9).Value = 25.02528114733049
$ws.Cells.Item(8, 11).Value = 8.047405245199259
$ws.Cells.Item(8, 12).Value = 10.14853174449693
$ws.Cells.Item(8, 13).Value = 13.38672358105632
$ws.Cells.Item(8, 15).Value = 24.53411164542691
$ws.Cells.Item(9, 2).Value = 11.26683569241673
$ws.Cells.Item(9, 3).Value = 9.556259359539879
$ws.Cells.Item(9, 4).Value = 6.249963555851145
$ws.Cells.Item(9, 5).Value = 12.76674512083656
$ws.Cells.Item(9, 7).Value = 33.69618503312819
$ws.Cells.Item(9, 8).Value = 15.57392067024147
$ws.Cells.Item(9, 9).Value = 24.81054603744043
$ws.Cells.Item(9, 11).Value = 8.542636793309521
$ws.Cells.Item(9, 12).Value = 10.13657947723817
$ws.Cells.Item(9, 13).Value = 13.4964462124835
$ws.Cells.Item(9, 15).Value = 24.35722545839207
$ws.Cells.Item(10, 2).Value = 11.72214722417597
$ws.Cells.Item(10, 3).Value = 9.601778924946929
$ws.Cells.Item(10, 4).Value = 6.463893772967144
$ws.Cells.Item(10, 5).Value = 12.72212788165073
$ws.Cells.Item(10, 7).Value = 33.60232375869845
$ws.Cells.Item(10, 8).Value = 15.5047219564316
$ws.Cells.Item(10, 9).Value = 24.67346170165577
$ws.Cells.Item(10, 11).Value = 8.892154333367783
$ws.Cells.Item(10, 12).Value = 10.13441882671696
$ws.Cells.Item(10, 13).Value = 13.58473485060826
$ws.Cells.Item(10, 15).Value = 24.25003987080702
$ws.Cells.Item(11, 2).Value = 11.92467332311905
$ws.Cells.Item(11, 3).Value = 9.622402155191653
$ws.Cells.Item(11, 4).Value = 6.560473236086461
$ws.Cells.Item(11, 5).Value = 12.70362070262339
$ws.Cells.Item(11, 7).Value = 33.56955380461488
$ws.Cells.Item(11, 8).Value = 15.47563443187629
$ws.Cells.Item(11, 9).Value = 24.6155907008053
$ws.Cells.Item(11, 11).Value = 9.047347617411287
$ws.Cells.Item(11, 12).Value = 10.13486382090444
$ws.Cells.Item(11, 13).Value = 13.62647075606706
$ws.Cells.Item(11, 15).Value = 24.20623893574522
$ws.Cells.Item(12, 2).Value = 12.00063298514309
$ws.Cells.Item(12, 3).Value = 9.630197821203113
$ws.Cells.Item(12, 4).Value = 6.596899751029496
$ws.Cells.Item(12, 5).Value = 12.69686947506925
$ws.Cells.Item(12, 7).Value = 33.55857684787498
$ws.Cells.Item(12, 8).Value = 15.46496343902323
$ws.Cells.Item(12, 9).Value = 24.59432228997646
$ws.Cells.Item(12, 11).Value = 9.105516059726796
$ws.Cells.Item(12, 12).Value = 10.13523668793862
$ws.Cells.Item(12, 13).Value = 13.64249266423032
$ws.Cells.Item(12, 15).Value = 24.19036690458671
$ws.Cells.Item(13, 2).Value = 11.98430736991542
$ws.Cells.Item(13, 3).Value = 9.628519532957098
$ws.Cells.Item(13, 4).Value = 6.589061755265938
$ws.Cells.Item(13, 5).Value = 12.69831204535349
$ws.Cells.Item(13, 7).Value = 33.56087716317725
$ws.Cells.Item(13, 8).Value = 15.46724633743656
$ws.Cells.Item(13, 9).Value = 24.59887407969504
$ws.Cells.Item(13, 11).Value = 9.093015919074009
$ws.Cells.Item(13, 12).Value = 10.13514731182259
$ws.Cells.Item(13, 13).Value = 13.63903253196421
$ws.Cells.Item(13, 15).Value = 24.19375343385417
$ws.Cells.Item(14, 2).Value = 11.93093764064093
$ws.Cells.Item(14, 3).Value = 9.623043798599799
$ws.Cells.Item(14, 4).Value = 6.56347321455568
$ws.Cells.Item(14, 5).Value = 12.70306012524257
$ws.Cells.Item(14, 7).Value = 33.56862199404774
$ws.Cells.Item(14, 8).Value = 15.4747496324399
$ws.Cells.Item(14, 9).Value = 24.61382798437837
$ws.Cells.Item(14, 11).Value = 9.05214547963929
$ws.Cells.Item(14, 12).Value = 10.13489040712152
$ws.Cells.Item(14, 13).Value = 13.62778458170696
$ws.Cells.Item(14, 15).Value = 24.20491880794664
$ws.Cells.Item(15, 2).Value = 11.89814971657659
$ws.Cells.Item(15, 3).Value = 9.619687887498548
$ws.Cells.Item(15, 4).Value = 6.547779331500187
$ws.Cells.Item(15, 5).Value = 12.70600192698167
$ws.Cells.Item(15, 7).Value = 33.57355258168372
$ws.Cells.Item(15, 8).Value = 15.47939039162454
$ws.Cells.Item(15, 9).Value = 24.62307183734192
$ws.Cells.Item(15, 11).Value = 9.027031540151253
$ws.Cells.Item(15, 12).Value = 10.13475962874624
$ws.Cells.Item(15, 13).Value = 13.62092293799342
$ws.Cells.Item(15, 15).Value = 24.21185100346289
$ws.Cells.Item(16, 2).Value = 11.70881291125839
$ws.Cells.Item(16, 3).Value = 9.600429303583565
$ws.Cells.Item(16, 4).Value = 6.457563763239273
$ws.Cells.Item(16, 5).Value = 12.72337335364038
$ws.Cells.Item(16, 7).Value = 33.60466554673752
$ws.Cells.Item(16, 8).Value = 15.50667100990107
$ws.Cells.Item(16, 9).Value = 24.67733408870834
$ws.Cells.Item(16, 11).Value = 8.881930917909271
$ws.Cells.Item(16, 12).Value = 10.1344183968669
$ws.Cells.Item(16, 13).Value = 13.58203822025929
$ws.Cells.Item(16, 15).Value = 24.25300226299565
$ws.Cells.Item(17, 2).Value = 11.59143067926069
$ws.Cells.Item(17, 3).Value = 9.58859200126277
$ws.Cells.Item(17, 4).Value = 6.402001270673526
$ws.Cells.Item(17, 5).Value = 12.73448828178788
$ws.Cells.Item(17, 7).Value = 33.62629871731468
$ws.Cells.Item(17, 8).Value = 15.52401922498223
$ws.Cells.Item(17, 9).Value = 24.7117724151039
$ws.Cells.Item(17, 11).Value = 8.791903031297062
$ws.Cells.Item(17, 12).Value = 10.13457411279271
$ws.Cells.Item(17, 13).Value = 13.55858043247896
$ws.Cells.Item(17, 15).Value = 24.27951821361841
$ws.Cells.Item(18, 2).Value = 11.52348669125588
$ws.Cells.Item(18, 3).Value = 9.581775717221415
$ws.Cells.Item(18, 4).Value = 6.369975938312216
$ws.Cells.Item(18, 5).Value = 12.74104974105211
$ws.Cells.Item(18, 7).Value = 33.63967578447289
$ws.Cells.Item(18, 8).Value = 15.53422254418245
$ws.Cells.Item(18, 9).Value = 24.73200295627091
$ws.Cells.Item(18, 11).Value = 8.739766192821076
$ws.Cells.Item(18, 12).Value = 10.13479813896405
$ws.Cells.Item(18, 13).Value = 13.54523673067968
$ws.Cells.Item(18, 15).Value = 24.29523610676557
$ws.Cells.Item(19, 2).Value = 11.50041064648918
$ws.Cells.Item(19, 3).Value = 9.57946656023012
$ws.Cells.Item(19, 4).Value = 6.359122408797297
$ws.Cells.Item(19, 5).Value = 12.74330027758192
$ws.Cells.Item(19, 7).Value = 33.64436533098557
$ws.Cells.Item(19, 8).Value = 15.53771587275448
$ws.Cells.Item(19, 9).Value = 24.73892522382552
$ws.Cells.Item(19, 11).Value = 8.722054238060093
$ws.Cells.Item(19, 12).Value = 10.13489711268802
$ws.Cells.Item(19, 13).Value = 13.54074456669628
$ws.Cells.Item(19, 15).Value = 24.30063801300617
$ws.Cells.Item(20, 2).Value = 11.60397113728186
$ws.Cells.Item(20, 3).Value = 9.589852920923876
$ws.Cells.Item(20, 4).Value = 6.407923249085143
$ws.Cells.Item(20, 5).Value = 12.73328764710895
$ws.Cells.Item(20, 7).Value = 33.6238991027281
$ws.Cells.Item(20, 8).Value = 15.52214918279021
$ws.Cells.Item(20, 9).Value = 24.70806266284489
$ws.Cells.Item(20, 11).Value = 8.8015238105819
$ws.Cells.Item(20, 12).Value = 10.13454362698016
$ws.Cells.Item(20, 13).Value = 13.56106223737278
$ws.Cells.Item(20, 15).Value = 24.27664724143407
$ws.Cells.Item(21, 2).Value = 11.94663404619573
$ws.Cells.Item(21, 3).Value = 9.624652546637989
$ws.Cells.Item(21, 4).Value = 6.570993459220969
$ws.Cells.Item(21, 5).Value = 12.70165852537163
$ws.Cells.Item(21, 7).Value = 33.56630824347764
$ws.Cells.Item(21, 8).Value = 15.47253640196326
$ws.Cells.Item(21, 9).Value = 24.6094181213997
$ws.Cells.Item(21, 11).Value = 9.0641667982178
$ws.Cells.Item(21, 12).Value = 10.13496032819704
$ws.Cells.Item(21, 13).Value = 13.63108254764946
$ws.Cells.Item(21, 15).Value = 24.20161986636775
$ws.Cells.Item(22, 2).Value = 12.16628502105826
$ws.Cells.Item(22, 3).Value = 9.647314206765445
$ws.Cells.Item(22, 4).Value = 6.676701163745363
$ws.Cells.Item(22, 5).Value = 12.68248515252204
$ws.Cells.Item(22, 7).Value = 33.53701952347284
$ws.Cells.Item(22, 8).Value = 15.44211561584381
$ws.Cells.Item(22, 9).Value = 24.54871432119487
$ws.Cells.Item(22, 11).Value = 9.232301202636052
$ws.Cells.Item(22, 12).Value = 10.13642331017281
$ws.Cells.Item(22, 13).Value = 13.67810794413598
$ws.Cells.Item(22, 15).Value = 24.15675001298257
$ws.Cells.Item(23, 2).Value = 12.04946862384286
$ws.Cells.Item(23, 3).Value = 9.635227360106033
$ws.Cells.Item(23, 4).Value = 6.620374998583327
$ws.Cells.Item(23, 5).Value = 12.69258136188
$ws.Cells.Item(23, 7).Value = 33.55188607418929
$ws.Cells.Item(23, 8).Value = 15.45816840759243
$ws.Cells.Item(23, 9).Value = 24.58076830625209
$ws.Cells.Item(23, 11).Value = 9.142902911767301
$ws.Cells.Item(23, 12).Value = 10.13553388684621
$ws.Cells.Item(23, 13).Value = 13.6528969978868
$ws.Cells.Item(23, 15).Value = 24.18031635701557
$ws.Cells.Item(24, 2).Value = 11.5983030215836
$ws.Cells.Item(24, 3).Value = 9.58928289284172
$ws.Cells.Item(24, 4).Value = 6.405246174126815
$ws.Cells.Item(24, 5).Value = 12.73382992032782
$ws.Cells.Item(24, 7).Value = 33.6249810412165
$ws.Cells.Item(24, 8).Value = 15.52299391366859
$ws.Cells.Item(24, 9).Value = 24.70973849781391
$ws.Cells.Item(24, 11).Value = 8.797175431481442
$ws.Cells.Item(24, 12).Value = 10.13455699058396
$ws.Cells.Item(24, 13).Value = 13.55993976886982
$ws.Cells.Item(24, 15).Value = 24.27794373307689
$ws.Cells.Item(25, 2).Value = 11.09559397678425
$ws.Cells.Item(25, 3).Value = 9.539434035611546
$ws.Cells.Item(25, 4).Value = 6.170802727565141
$ws.Cells.Item(25, 5).Value = 12.78476430059899
$ws.Cells.Item(25, 7).Value = 33.73951749984511
$ws.Cells.Item(25, 8).Value = 15.60152447273358
$ws.Cells.Item(25, 9).Value = 24.86500606859686
$ws.Cells.Item(25, 11).Value = 8.542636793309521
$ws.Cells.Item(25, 12).Value = 10.13864712299851
$ws.Cells.Item(25, 13).Value = 13.4964462124835
$ws.Cells.Item(25, 15).Value = 24.40108414001183
